$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.118.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.63%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.857.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.61%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.39"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.83%  "

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.07%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.59%  "

# Row 8
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.84"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.24%  "

# Row 9
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2828"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.57%  "

# Row 10
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06453"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.96%  "

# Row 11
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.96"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.53%  "

# Row 12
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07730"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.67%  "

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.864.58"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.29%  "

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.52"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.63%  "

# Row 15
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.052"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.31%  "

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.6782"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.03%  "

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "265.60"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.09%  "

# Row 18
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "30.091.60"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.79%  "

# Row 19
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.33"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.54%  "

# Row 20
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007544"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.52%  "

# Row 21
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.04%  "

# Row 22
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.106.45"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.43%  "

# Row 23
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.02%  "

# Row 24
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.146"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.50%  "

# Row 25
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.087"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.00%  "

# Row 26
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.270"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.55%  "

# Row 27
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.09"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.16%  "

# Row 28
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.48"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.16%  "

# Row 29
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.881"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.63%  "

# Row 30
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.365"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.42%  "

# Row 31
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09832"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.34%  "

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.451"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.69%  "

# Row 33
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.199"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.02%  "

# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.970"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.45%  "

# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.04651"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.03%  "

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.112"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.98%  "

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6842"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.42%  "

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.712"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.33%  "

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01807"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.60%  "

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.716"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.37%  "

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.277"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.21%  "

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.45"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.63%  "

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.03%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8308"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.29%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.875"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.20%  "

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.99"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.97%  "

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "TheSandbox"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4036"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.13%  "

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.120"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.51%  "

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "926.19"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.03%  "

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.924"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.97%  "

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.07"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.32%  "
